$wb = $excel.ActiveWorkbook
$n = $wb.Names.Add("TestName", "=1+1")
try {
  $n.Comment = "hello comment"
  Write-Host "comment set ok"
} catch {
  Write-Host "comment set failed: $_"
}
Write-Host "Comment prop value:" $n.Comment
